$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, pushing existing rows 52:88 down to 53:89
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record's data
$ws.Range("A52").Value = 4
$ws.Range("B52").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C52").Value = "Los Lagos"
$ws.Range("D52").Value = 44488
$ws.Range("E52").Value = 10
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100108
$ws.Range("H52").Value = "Tropicales y subtropicales"
$ws.Range("I52").Value = 100108002
$ws.Range("J52").Value = "Mango"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 200
$ws.Range("N52").Value = 7500
$ws.Range("O52").Value = 8000
$ws.Range("P52").Value = 7750
$ws.Range("Q52").Value = "$/bandeja 4 kilos"
$ws.Range("R52").Value = "Perú"
$ws.Range("S52").Value = 1938
$ws.Range("T52").Value = 4
